# Update "想去人数" (want-to-go count) figures for a handful of events.
# These numbers appear both on the "展览" sheet and on the combined
# "全部类型" sheet, which lists the same events (shifted by one row
# because it also contains an extra "演出" entry).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 31
$ws1.Range("F8").Value = 8058
$ws1.Range("F12").Value = 789
$ws1.Range("F13").Value = 37
$ws1.Range("F15").Value = 203
$ws1.Range("F16").Value = 60
$ws1.Range("F17").Value = 52
$ws1.Range("F19").Value = 857

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 31
$ws4.Range("F9").Value = 8058
$ws4.Range("F13").Value = 789
$ws4.Range("F14").Value = 37
$ws4.Range("F16").Value = 203
$ws4.Range("F17").Value = 60
$ws4.Range("F18").Value = 52
$ws4.Range("F20").Value = 857
